$d = $word.ActiveDocument

# Original paragraph text is "Version 1." (positions 0..10):
#   V  e  r  s  i  o  n  sp 1  .
#   0  1  2  3  4  5  6  7  8  9 (end=10)
#
# Target paragraph text is "Version 2." but the run layout differs from a
# simple text substitution:
#   - "Version" is split into two runs: "Versi" + "on"
#   - " 1." becomes " 2" (no trailing period)
#   - the existing "_GoBack" bookmark ends up between " 2" and the final "."
#   - "." becomes its own trailing run, placed after the bookmark

# Step 1: split "Version" into "Versi" | "on" (clean run break, no format
# residue) by adding then immediately removing a bookmark at the boundary.
$b1 = $d.Bookmarks.Add("TmpSplit1", $d.Range(5, 5))
$d.Bookmarks("TmpSplit1").Delete()

# Step 2: change the "1" digit to "2" in place (stays inside the " 1." run).
$d.Range(8, 9).Text = "2"

# Step 3: drop the trailing "." so the run becomes " 2", then retype the
# "." afterwards so it lands in a fresh run at the very end of the story,
# after the "_GoBack" bookmark which already sits there.
$d.Range(9, 10).Delete()
$d.Range(9, 9).InsertAfter(".")
